$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","I","J","K","M","O")

$data = @{
    2 = @(4.342782668993784, 9.055350555275817, 13.32567920946257, 32.32900490095896, 3.655543838190128, 22.4992872989442, 9.918680427685027, 14.87214931442036, 17.34426775029964, 24.41475935850357)
    3 = @(4.231923666515726, 9.030910837732627, 13.35114581393254, 32.47284114800821, 3.657686912343961, 22.60621924636345, 9.951922749143936, 14.32033050665856, 17.11510346478222, 24.55090742061739)
    4 = @(4.161722679053244, 9.016977391424312, 13.36904847266762, 32.57051913536, 3.65907176536805, 22.67817167450871, 9.97357989110702, 13.97065735200205, 16.97447443289569, 24.64110421674761)
    5 = @(4.132604443711554, 9.011572637222011, 13.37691351846389, 32.61266854613606, 3.659653509387644, 22.70907038010944, 9.982719224268907, 13.82562219586961, 16.91724364144145, 24.67951547838259)
    6 = @(4.127739319064112, 9.010691783472884, 13.37825389493055, 32.61980877042217, 3.659751160396421, 22.71429619868649, 9.98425577546495, 13.80139170513689, 16.90774680808229, 24.68599348309128)
    7 = @(4.161332013203991, 9.016903390190995, 13.36915223765828, 32.57107809536765, 3.659079540425173, 22.67858200549102, 9.973701875830463, 13.96871137880696, 16.97370221272142, 24.64161554817311)
    8 = @(4.305016011907021, 9.046703281283559, 13.33398964292108, 32.37665103556707, 3.656268483258158, 22.53484738938645, 9.929884027044858, 14.68424193866269, 17.26526960571447, 24.46033031326352)
    9 = @(4.56882386230467, 9.113489626133106, 13.2830253161047, 32.07007006584518, 3.651300936316999, 22.30318655638213, 9.85382286163404, 15.99399766303063, 17.83517820200335, 24.1574170895962)
    10 = @(4.750448452672321, 9.167419248266434, 13.2565581668719, 31.89090096418276, 3.647979931170894, 22.16392464154621, 9.80392333069975, 16.89093097512941, 18.24948964543069, 23.96721140565472)
    11 = @(4.830198203550448, 9.192956481343709, 13.2469013818013, 31.81950844825057, 3.646539727229774, 22.10735896258548, 9.782515209712727, 17.28335116128769, 18.4363428190357, 23.88776796178604)
    12 = @(4.859967568113243, 9.202766448363155, 13.24358717713646, 31.79393685169783, 3.646004445470666, 22.08692038106853, 9.774593745749184, 17.42960857283956, 18.50681210861852, 23.85870857948654)
    13 = @(4.853575588382496, 9.200647568065182, 13.24428571534054, 31.79937893584763, 3.646119279930096, 22.09127844515516, 9.776291538205161, 17.39821510951805, 18.49164895124977, 23.86492139732896)
    14 = @(4.832656059874521, 9.193760784426136, 13.24662185468149, 31.81737528227712, 3.646495487305825, 22.10565776183906, 9.781859794011586, 17.29543134265978, 18.44214645014952, 23.88535667090215)
    15 = @(4.819785748261056, 9.189560459906634, 13.24809742085625, 31.82858937872674, 3.646727237981558, 22.1145935165281, 9.785294635435573, 17.23216536450089, 18.41178563401865, 23.89800740839281)
    16 = @(4.745177322152114, 9.165770129270259, 13.25723721142916, 31.89577090147817, 3.648075466916739, 22.167758385799, 9.805348355942661, 16.86496259405465, 18.23724101707561, 23.97254621982644)
    17 = @(4.698659127059138, 9.151429374645048, 13.26345454594305, 31.93958116151886, 3.648920591717389, 22.20211541736805, 9.817981156431641, 16.6356243377365, 18.12970906631754, 24.02009131783329)
    18 = @(4.671633663515443, 9.143275790248758, 13.26725492745769, 31.96573097084933, 3.649413327543781, 22.2225149888878, 9.825368786242651, 16.50225183473304, 18.0677103186837, 24.04810425294227)
    19 = @(4.662437580098352, 9.140531560497459, 13.26858020253577, 31.97474795375648, 3.649581301827598, 22.22953138285667, 9.827891003454663, 16.45684617846388, 18.04669468868314, 24.05770320714679)
    20 = @(4.703639082355337, 9.152946191705071, 13.26276948263166, 31.93481896126598, 3.648829939646734, 22.19839194626572, 9.816623793251436, 16.66018999391872, 18.14117187399667, 24.01496106498394)
    21 = @(4.838812434550051, 9.195779851323428, 13.24592637685575, 31.81204953782381, 3.64638471258495, 22.10140752086037, 9.780219235553941, 17.32568576401511, 18.45669477988297, 23.87932649528634)
    22 = @(4.924642073629655, 9.224585284761842, 13.23691531314081, 31.74034522639818, 3.644845417079804, 22.0437468327623, 9.757506810886781, 17.7469281791173, 18.6612045834224, 23.79665354582124)
    23 = @(4.879068443695908, 9.209138708004083, 13.24154203077851, 31.77783153170962, 3.645661604471125, 22.07399574256368, 9.769530156408459, 17.52338622717629, 18.55222737768672, 23.84022930244633)
    24 = @(4.701388520175209, 9.152260154573014, 13.26307849592755, 31.93696895387293, 3.648870902067479, 22.2000733118015, 9.817237068165106, 16.64908859645415, 18.13599008635607, 24.01727833880592)
    25 = @(4.499515159442733, 9.094551103817244, 13.2948853014343, 32.14495522574269, 3.652586819563192, 22.36044755776615, 9.873346562932465, 15.65061046354425, 17.68154269250163, 24.23370403585931)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
